$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '27.002.05'
Set-TextCell $ws.Range('E2') '  -3.21%  '
Set-TextCell $ws.Range('D3') '1.718.17'
Set-TextCell $ws.Range('E3') '  -3.03%  '
Set-TextCell $ws.Range('D4') '1.012'
Set-TextCell $ws.Range('E4') '  +0.80%  '
Set-TextCell $ws.Range('D5') '319.38'
Set-TextCell $ws.Range('E5') '  -2.47%  '
Set-TextCell $ws.Range('D6') '1.011'
Set-TextCell $ws.Range('E6') '  +0.94%  '
Set-TextCell $ws.Range('D7') '0.4651'
Set-TextCell $ws.Range('E7') '  +3.59%  '
Set-TextCell $ws.Range('D8') '0.3440'
Set-TextCell $ws.Range('E8') '  -3.74%  '
Set-TextCell $ws.Range('D9') '42.00'
Set-TextCell $ws.Range('E9') '  -0.47%  '
Set-TextCell $ws.Range('D10') '0.07276'
Set-TextCell $ws.Range('E10') '  -2.58%  '
Set-TextCell $ws.Range('D11') '1.049'
Set-TextCell $ws.Range('E11') '  -4.34%  '
Set-TextCell $ws.Range('D12') '1.013'
Set-TextCell $ws.Range('E12') '  +1.18%  '
Set-TextCell $ws.Range('D13') '19.80'
Set-TextCell $ws.Range('E13') '  -5.31%  '
Set-TextCell $ws.Range('D14') '5.856'
Set-TextCell $ws.Range('E14') '  -3.30%  '
Set-TextCell $ws.Range('D15') '1.723.73'
Set-TextCell $ws.Range('E15') '  -2.85%  '
Set-TextCell $ws.Range('D16') '6.886'
Set-TextCell $ws.Range('E16') '  -4.64%  '
Set-TextCell $ws.Range('D17') '89.80'
Set-TextCell $ws.Range('E17') '  -3.54%  '
Set-TextCell $ws.Range('D18') '0.00001042'
Set-TextCell $ws.Range('E18') '  -1.76%  '
Set-TextCell $ws.Range('D19') '0.06306'
Set-TextCell $ws.Range('E19') '  -1.92%  '
Set-TextCell $ws.Range('D20') '1.010'
Set-TextCell $ws.Range('E20') '  +0.84%  '
Set-TextCell $ws.Range('D21') '16.41'
Set-TextCell $ws.Range('E21') '  -4.98%  '
Set-TextCell $ws.Range('D22') '5.608'
Set-TextCell $ws.Range('E22') '  -3.79%  '
Set-TextCell $ws.Range('D23') '27.076.37'
Set-TextCell $ws.Range('E23') '  -3.08%  '
Set-TextCell $ws.Range('D24') '10.79'
Set-TextCell $ws.Range('E24') '  -4.93%  '
Set-TextCell $ws.Range('D25') '2.110'
Set-TextCell $ws.Range('E25') '  -0.22%  '
Set-TextCell $ws.Range('D26') '157.34'
Set-TextCell $ws.Range('E26') '  -3.18%  '
Set-TextCell $ws.Range('D27') '19.47'
Set-TextCell $ws.Range('E27') '  -4.00%  '
Set-TextCell $ws.Range('D28') '1.933.34'
Set-TextCell $ws.Range('E28') '  -2.19%  '
Set-TextCell $ws.Range('D29') '2.102'
Set-TextCell $ws.Range('E29') '  -4.56%  '
Set-TextCell $ws.Range('D30') '119.23'
Set-TextCell $ws.Range('E30') '  -5.24%  '
Set-TextCell $ws.Range('D31') '1.014'
Set-TextCell $ws.Range('E31') '  -8.26%  '
Set-TextCell $ws.Range('D32') '0.09097'
Set-TextCell $ws.Range('E32') '  -0.79%  '
Set-TextCell $ws.Range('D33') '3.583'
Set-TextCell $ws.Range('E33') '  -1.36%  '
Set-TextCell $ws.Range('D34') '5.312'
Set-TextCell $ws.Range('E34') '  -4.69%  '
Set-TextCell $ws.Range('D35') '0.02189'
Set-TextCell $ws.Range('E35') '  -4.91%  '
Set-TextCell $ws.Range('D36') '11.11'
Set-TextCell $ws.Range('E36') '  -6.59%  '
Set-TextCell $ws.Range('D37') '0.05798'
Set-TextCell $ws.Range('E37') '  -5.32%  '
Set-TextCell $ws.Range('D38') '0.1993'
Set-TextCell $ws.Range('E38') '  -4.97%  '
Set-TextCell $ws.Range('B39') 'WEMIXTOKEN'
Set-TextCell $ws.Range('C39') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws.Range('D39') '1.413'
Set-TextCell $ws.Range('E39') '  +1.60%  '
Set-TextCell $ws.Range('B40') 'InternetComputer(DFINITY)'
Set-TextCell $ws.Range('C40') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws.Range('D40') '4.744'
Set-TextCell $ws.Range('E40') '  -4.71%  '
Set-TextCell $ws.Range('D41') '0.5929'
Set-TextCell $ws.Range('E41') '  -6.70%  '
Set-TextCell $ws.Range('D42') '1.133'
Set-TextCell $ws.Range('E42') '  -4.43%  '
Set-TextCell $ws.Range('D43') '7.506'
Set-TextCell $ws.Range('E43') '  -5.59%  '
Set-TextCell $ws.Range('D44') '3.656'
Set-TextCell $ws.Range('E44') '  -2.24%  '
Set-TextCell $ws.Range('D45') '12.43'
Set-TextCell $ws.Range('E45') '  -5.80%  '
Set-TextCell $ws.Range('D46') '0.5560'
Set-TextCell $ws.Range('E46') '  -5.50%  '
Set-TextCell $ws.Range('D47') '120.13'
Set-TextCell $ws.Range('E47') '  -1.96%  '
Set-TextCell $ws.Range('D48') '1.856'
Set-TextCell $ws.Range('E48') '  -5.44%  '
Set-TextCell $ws.Range('B49') 'Cronos'
Set-TextCell $ws.Range('C49') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D49') '0.06683'
Set-TextCell $ws.Range('E49') '  -3.45%  '
Set-TextCell $ws.Range('B50') 'EOS'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextCell $ws.Range('D50') '1.088'
Set-TextCell $ws.Range('E50') '  -4.56%  '
Set-TextCell $ws.Range('B51') 'PaxDollar'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws.Range('D51') '1.012'
Set-TextCell $ws.Range('E51') '  +1.08%  '
